# "Added EventId in Second Position"
#
# The header row currently reads:
#   A1 CoID | B1 Exhibiting As | C1 Booth Number | D1 Company Contact First Name |
#   E1 ... | ... | L1 EventId
#
# EventId needs to move out of the last column (L) into the second column (B),
# pushing every header between B and L one slot to the right so the rest of
# the header order is preserved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rotate the header labels (B1..L1) so EventId lands in B1 -----------
# Snapshot the current labels for columns B (2) through L (12) left-to-right.
$oldHeaders = @()
for ($col = 2; $col -le 12; $col++) {
    $oldHeaders += $ws.Cells.Item(1, $col).Value2
}

# Last entry (column L / index 10) is "EventId" - move it to the front, and
# shift the remaining ten headers one column to the right.
$eventId   = $oldHeaders[10]
$rest      = $oldHeaders[0..9]
$newHeaders = @($eventId) + $rest

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $col = 2 + $i
    $ws.Cells.Item(1, $col).Value = $newHeaders[$i]
}

# --- 2. Fix up header formatting that travelled with the old labels --------
# "Exhibiting As" (now in C1) had inherited "Booth Number"'s centered look;
# restore the normal (non-centered) header formatting used by the other
# columns.
$ws.Range("C1").HorizontalAlignment = 1

# "Booth Number" (now in D1) used to be plain; give it the centered header
# style that used to live on the old Booth Number column (now C1).
$ws.Range("D1").HorizontalAlignment = -4108

# New EventId header (B1) keeps the same bold Calibri / bordered / wrapped
# header look shared by the rest of row 1 - make that explicit.
$b1 = $ws.Range("B1")
$b1.Font.Bold = $true
$b1.Font.Name = "Calibri"
$b1.Font.Size = 10
$b1.Font.Color = 0
$b1.WrapText = $true
$b1.Borders.LineStyle = 1
$b1.HorizontalAlignment = 1

# --- 3. Row height grew to fit the new header text --------------------------
$ws.Rows(1).RowHeight = 35.25

# --- 4. Selection moves to the newly inserted EventId header ---------------
$null = $ws.Range("B1").Select()
